# Updates cryptos list price/volume(1h) figures on the Sheet1 table.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h)
# Note: some "Price" values (e.g. 0.0794, 0.130, 0.0958) look like plain
# numbers to Excel's auto-detection and would otherwise be stored as a
# number (losing significant trailing zeros). Prefixing the literal value
# with a leading apostrophe forces Excel to keep them as text, matching
# the original inline-string cell content exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.930.67'
$ws.Range('E2').Value = '  -1.08%  '
$ws.Range('D3').Value = '1.636.99'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').Value = "'215.55"
$ws.Range('E5').Value = '  -0.61%  '
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('E8').Value = '  -0.96%  '
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('E10').Value = '  -1.62%  '
$ws.Range('D11').Value = "'0.0794"
$ws.Range('E11').Value = '  +0.01%  '
$ws.Range('D12').Value = '1.863.85'
$ws.Range('E12').Value = '  -0.44%  '
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('D14').Value = '1.651.51'
$ws.Range('E14').Value = '  +1.06%  '
$ws.Range('D15').Value = "'0.544"
$ws.Range('E15').Value = '  -0.54%  '
$ws.Range('E16').Value = '  -0.31%  '
$ws.Range('E17').Value = '  -0.62%  '
$ws.Range('D18').Value = '25.907.67'
$ws.Range('E18').Value = '  -1.17%  '
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('D20').Value = "'192.74"
$ws.Range('E20').Value = '  -1.17%  '
$ws.Range('E21').Value = '  -2.20%  '
$ws.Range('E22').Value = '  -1.50%  '
$ws.Range('E23').Value = '  -0.84%  '
$ws.Range('D24').Value = "'0.130"
$ws.Range('E24').Value = '  +4.76%  '
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('D26').Value = "'143.42"
$ws.Range('E26').Value = '  +0.25%  '
$ws.Range('E27').Value = '  +0.22%  '
$ws.Range('E28').Value = '  -0.98%  '
$ws.Range('E29').Value = '  -0.40%  '
$ws.Range('E30').Value = '  -0.35%  '
$ws.Range('E31').Value = '  -0.33%  '
$ws.Range('E32').Value = '  -1.71%  '
$ws.Range('E33').Value = '  -0.22%  '
$ws.Range('E35').Value = '  +1.68%  '
$ws.Range('E36').Value = '  -1.21%  '
$ws.Range('D37').Value = '1.131.18'
$ws.Range('E37').Value = '  -0.45%  '
$ws.Range('E38').Value = '  -1.76%  '
$ws.Range('E39').Value = '  -1.28%  '
$ws.Range('E40').Value = '  -0.41%  '
$ws.Range('E41').Value = '  -0.82%  '
$ws.Range('D42').Value = "'99.27"
$ws.Range('E42').Value = '  -1.06%  '
$ws.Range('D43').Value = "'0.796"
$ws.Range('E43').Value = '  -0.50%  '
$ws.Range('D44').Value = '1.773.74'
$ws.Range('E44').Value = '  -0.43%  '
$ws.Range('E45').Value = '  +3.27%  '
$ws.Range('D46').Value = "'56.61"
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('E47').Value = '  +2.30%  '
$ws.Range('E48').Value = '  -0.57%  '
$ws.Range('D49').Value = "'7.67"
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('E50').Value = '  -0.83%  '
$ws.Range('D51').Value = "'0.0958"
$ws.Range('E51').Value = '  -1.16%  '
